$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that looks numeric (e.g. "0.9994", "329.47").
# Force text format first so Excel does not coerce these into numbers,
# matching the original workbook where these are stored as text.
$dCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D19","D20","D22","D23","D24","D25","D26","D27","D28","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '28.337.83'
$ws.Range("E2").Value = '  +4.29%  '
$ws.Range("D3").Value = '1.818.11'
$ws.Range("E3").Value = '  +4.06%  '
$ws.Range("D4").Value = '0.9994'
$ws.Range("E4").Value = '  -1.75%  '
$ws.Range("D5").Value = '329.47'
$ws.Range("E5").Value = '  +2.05%  '
$ws.Range("D6").Value = '0.9967'
$ws.Range("E6").Value = '  -1.68%  '
$ws.Range("D7").Value = '0.4434'
$ws.Range("E7").Value = '  +5.71%  '
$ws.Range("D8").Value = '0.3742'
$ws.Range("E8").Value = '  +4.94%  '
$ws.Range("D9").Value = '44.77'
$ws.Range("E9").Value = '  -0.12%  '
$ws.Range("D10").Value = '0.07704'
$ws.Range("E10").Value = '  +5.03%  '
$ws.Range("D11").Value = '1.127'
$ws.Range("E11").Value = '  +1.72%  '
$ws.Range("D12").Value = '0.9969'
$ws.Range("E12").Value = '  -1.99%  '
$ws.Range("D13").Value = '22.06'
$ws.Range("E13").Value = '  +2.98%  '
$ws.Range("D14").Value = '6.316'
$ws.Range("E14").Value = '  +4.27%  '
$ws.Range("D15").Value = '7.505'
$ws.Range("E15").Value = '  +4.75%  '
$ws.Range("D16").Value = '1.820.78'
$ws.Range("E16").Value = '  +3.93%  '
$ws.Range("D17").Value = '93.67'
$ws.Range("E17").Value = '  +11.66%  '
$ws.Range("D19").Value = '0.06489'
$ws.Range("E19").Value = '  +11.29%  '
$ws.Range("D20").Value = '0.9987'
$ws.Range("E20").Value = '  -1.37%  '
$ws.Range("E21").Value = '  +5.04%  '
$ws.Range("D22").Value = '6.256'
$ws.Range("E22").Value = '  +3.24%  '
$ws.Range("D23").Value = '0.5346'
$ws.Range("E23").Value = '  -1.92%  '
$ws.Range("D24").Value = '28.392.23'
$ws.Range("E24").Value = '  +4.20%  '
$ws.Range("D25").Value = '11.74'
$ws.Range("E25").Value = '  +5.04%  '
$ws.Range("D26").Value = '2.139'
$ws.Range("E26").Value = '  -11.78%  '
$ws.Range("D27").Value = '20.60'
$ws.Range("E27").Value = '  +3.93%  '
$ws.Range("D28").Value = '155.37'
$ws.Range("E28").Value = '  +4.57%  '
$ws.Range("E29").Value = '  +1.68%  '
$ws.Range("D30").Value = '2.022.24'
$ws.Range("E30").Value = '  +3.47%  '
$ws.Range("D31").Value = '127.67'
$ws.Range("E31").Value = '  +1.15%  '
$ws.Range("D32").Value = '1.200'
$ws.Range("E32").Value = '  -2.58%  '
$ws.Range("D33").Value = '5.849'
$ws.Range("E33").Value = '  +6.34%  '
$ws.Range("D34").Value = '0.09228'
$ws.Range("E34").Value = '  +2.50%  '
$ws.Range("D35").Value = '3.675'
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("D36").Value = '13.04'
$ws.Range("E36").Value = '  +5.33%  '
$ws.Range("D37").Value = '0.02341'
$ws.Range("E37").Value = '  +4.76%  '
$ws.Range("D38").Value = '0.2173'
$ws.Range("E38").Value = '  +1.87%  '
$ws.Range("D39").Value = '5.172'
$ws.Range("E39").Value = '  +4.48%  '
$ws.Range("D40").Value = '0.6570'
$ws.Range("E40").Value = '  +2.92%  '
$ws.Range("D41").Value = '0.06196'
$ws.Range("E41").Value = '  +2.45%  '
$ws.Range("D42").Value = '1.195'
$ws.Range("E42").Value = '  +1.99%  '
$ws.Range("D43").Value = '8.077'
$ws.Range("E43").Value = '  +3.03%  '
$ws.Range("D44").Value = '0.9958'
$ws.Range("E44").Value = '  -1.73%  '
$ws.Range("D45").Value = '14.02'
$ws.Range("E45").Value = '  +4.10%  '
$ws.Range("D46").Value = '1.392'
$ws.Range("E46").Value = '  -1.73%  '
$ws.Range("D47").Value = '0.6077'
$ws.Range("E47").Value = '  +4.62%  '
$ws.Range("D48").Value = '3.760'
$ws.Range("E48").Value = '  +0.30%  '
$ws.Range("D49").Value = '126.80'
$ws.Range("E49").Value = '  +3.23%  '
$ws.Range("D50").Value = '2.036'
$ws.Range("D51").Value = '0.06989'
$ws.Range("E51").Value = '  +2.34%  '

Write-Host "Applied crypto price/volume updates."